$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark from the Cluster 8 block ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Append a duplicate of the Cluster 8 parameter block as "Cluster 9",
#        with updated max_pod_connected / min_pod_connected values ---
$lines = @(
    "Cluster 9",
    "num_pods=100",
    "max_pod_connected=25",
    "min_pod_connected=5",
    "num_cores=5",
    "num_slots=80",
    "",
    "data_rate_choice = [50, 400]",
    "",
    "data_rate_probs = {}",
    "data_rate_probs['90_10'] = [0.90, 0.10]"
)

$prevCount = $d.Paragraphs.Count
foreach ($line in $lines) {
    $cur = $d.Paragraphs.Last
    $cur.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    if ($line -ne "") {
        $newPara.Range.InsertAfter($line)
    }
}

# The first of the newly-added paragraphs is the new "Cluster 9" heading.
$headingIndex = $prevCount + 1
$headingPara = $d.Paragraphs.Item($headingIndex)
$headingPara.Style = "Heading 1"

# Re-create the "_GoBack" bookmark inside the new heading, between "C" and "luster 9"
$bmPos = $headingPara.Range.Start + 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
